# Generate Report for Handoff
#
# - "Handed back: in sync with en-US" -> "Ready for handoff"
#     (Overview!E2, Overview!F2, zh-cn!C2, de-de!C2 — "Status" columns)
# - Latest HO Xliff Generate Date / Latest Handoff Datetime timestamps bumped
#     (Overview!G2: 05:04:09 -> 05:04:57, de-de!H2: 05:04:09 -> 05:04:57,
#      zh-cn!H2: 05:04:03 -> 05:04:53)
# - Status columns on the per-language sheets (and mirrored on Overview)
#     narrowed from ~29.98 chars wide to ~17.22 chars wide

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Status text: "Handed back: in sync with en-US" -> "Ready for handoff" ---
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("C2").Value = "Ready for handoff"

# --- Timestamps ---
$wsOverview.Range("G2").Value = "2016-08-21 05:04:57"
$wsDeDe.Range("H2").Value = "2016-08-21 05:04:57"
$wsZhCn.Range("H2").Value = "2016-08-21 05:04:53"

# --- Column widths: stored width 29.9777047293527 -> 17.2159881591797 ---
# ColumnWidth is expressed in characters; the stored <col width> includes the
# default ~0.8333 char padding, so back that out of the target before setting.
$targetColumnWidth = 17.2159881591797 - 0.8333333333333333

$wsOverview.Range("E1:F1").ColumnWidth = $targetColumnWidth
$wsZhCn.Range("C1").ColumnWidth = $targetColumnWidth
$wsDeDe.Range("C1").ColumnWidth = $targetColumnWidth
